# Disaggregation of commodity Copper
#
# 1) Rename the shared label "Copper ores and concentrates" -> "Copper".
#    The label lives in cell C4 of every per-year worksheet (2000..2100)
#    and all of them point at the same shared string, so updating it on
#    any one sheet (via the cell text) updates it everywhere.
# 2) A handful of years got refreshed D4 totals (last-digit precision
#    updates coming from the disaggregated source data).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Cells.Item(4, 3)
    if ($cell.Value2 -eq "Copper ores and concentrates") {
        $cell.Value2 = "Copper"
    }
}

$updates = @{
    "2033" = 95388.41488819558
    "2047" = 634503.609349301
    "2048" = 806653.9762728701
    "2054" = 1998651.653451595
    "2065" = 909749.2385804425
    "2073" = 879339.2488812557
}

foreach ($year in $updates.Keys) {
    $sheetName = [string]$year
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Item(4, 4).Value2 = $updates[$year]
}
